$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel
# auto-converting numeric-looking strings (e.g. "1.013", "27.486.29")
# into actual numbers. We temporarily force a text number-format,
# assign the value, then restore the cell's original style so the
# resulting cell keeps its original formatting/style index.
function Set-TextValue($cellRef, [string]$text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-TextValue "D2" '27.486.29'
Set-TextValue "E2" '  +1.91%  '
Set-TextValue "D3" '1.866.28'
Set-TextValue "E3" '  +0.95%  '
Set-TextValue "D4" '1.013'
Set-TextValue "E4" '  -0.15%  '
Set-TextValue "D5" '311.78'
Set-TextValue "E5" '  +0.71%  '
Set-TextValue "E6" '  +0.05%  '
Set-TextValue "D7" '0.4776'
Set-TextValue "E7" '  +0.22%  '
Set-TextValue "D8" '0.3775'
Set-TextValue "E8" '  +2.84%  '
Set-TextValue "D9" '0.07339'
Set-TextValue "E9" '  +1.64%  '
Set-TextValue "D10" '0.9359'
Set-TextValue "E10" '  +0.98%  '
Set-TextValue "D11" '20.68'
Set-TextValue "E11" '  +5.03%  '
Set-TextValue "D12" '0.07833'
Set-TextValue "E12" '  +1.55%  '
Set-TextValue "D13" '1.884.49'
Set-TextValue "E13" '  +2.91%  '
Set-TextValue "D14" '5.438'
Set-TextValue "E14" '  +2.20%  '
Set-TextValue "D15" '6.553'
Set-TextValue "E15" '  +2.17%  '
Set-TextValue "D16" '90.50'
Set-TextValue "E16" '  +1.87%  '
Set-TextValue "E17" '  -0.19%  '
Set-TextValue "D18" '0.000008896'
Set-TextValue "E18" '  +3.03%  '
Set-TextValue "E19" '  -0.06%  '
Set-TextValue "D20" '27.557.81'
Set-TextValue "E20" '  +2.06%  '
Set-TextValue "D21" '14.73'
Set-TextValue "E21" '  +1.27%  '
Set-TextValue "D22" '5.120'
Set-TextValue "E22" '  +1.23%  '
Set-TextValue "D23" '10.70'
Set-TextValue "E23" '  +0.37%  '
Set-TextValue "D24" '1.936'
Set-TextValue "E24" '  +0.13%  '
Set-TextValue "D25" '154.58'
Set-TextValue "E25" '  +1.31%  '
Set-TextValue "D26" '18.48'
Set-TextValue "E26" '  +1.54%  '
Set-TextValue "D27" '2.025'
Set-TextValue "E27" '  +1.42%  '
Set-TextValue "D28" '115.52'
Set-TextValue "E28" '  +1.16%  '
Set-TextValue "E29" '  +0.56%  '
Set-TextValue "D30" '0.08900'
Set-TextValue "E30" '  +0.29%  '
Set-TextValue "D31" '3.330'
Set-TextValue "E31" '  +0.30%  '
Set-TextValue "E32" '  +4.03%  '
Set-TextValue "D33" '0.7583'
Set-TextValue "E33" '  +2.00%  '
Set-TextValue "D34" '4.612'
Set-TextValue "E34" '  +2.59%  '
Set-TextValue "D35" '2.747'
Set-TextValue "E35" '  +1.21%  '
Set-TextValue "D36" '1.120'
Set-TextValue "E36" '  +0.66%  '
Set-TextValue "D37" '0.02035'
Set-TextValue "E37" '  +4.08%  '
Set-TextValue "D38" '0.05266'
Set-TextValue "E38" '  +0.05%  '
Set-TextValue "D39" '2.992'
Set-TextValue "E39" '  +0.36%  '
Set-TextValue "D40" '0.5318'
Set-TextValue "E40" '  +2.48%  '
Set-TextValue "D41" '7.087'
Set-TextValue "E41" '  +1.45%  '
Set-TextValue "B42" 'Aptos'
Set-TextValue "C42" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D42" '8.499'
Set-TextValue "E42" '  +3.75%  '
Set-TextValue "B43" 'Algorand'
Set-TextValue "C43" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D43" '0.1524'
Set-TextValue "E43" '  +1.00%  '
Set-TextValue "D44" '10.62'
Set-TextValue "E44" '  -0.01%  '
Set-TextValue "D45" '0.4806'
Set-TextValue "E45" '  +1.71%  '
Set-TextValue "E46" '  +0.06%  '
Set-TextValue "D47" '102.99'
Set-TextValue "E47" '  +1.67%  '
Set-TextValue "D48" '1.652'
Set-TextValue "E48" '  +3.10%  '
Set-TextValue "D49" '67.47'
Set-TextValue "E49" '  +3.03%  '
Set-TextValue "D50" '0.06080'
Set-TextValue "E50" '  +0.85%  '
Set-TextValue "D51" '0.9195'
Set-TextValue "E51" '  +3.56%  '
Write-Host "Applied cryptos list update."
